# [F2] Add filename check for assignment
#
# For every "success" row, the "Expected filename on server" column (E) is
# updated from the bare uploaded name to the filename actually expected to
# be stored on the server (an extension appended / path separators
# stripped). "fail" rows keep their "None" expectation untouched.

$wb = $excel.ActiveWorkbook

# Sheet F2-1
$ws1 = $wb.Worksheets.Item("F2-1")
$ws1.Range("E2").Value = "testing1.php"
$ws1.Range("E3").Value = ($ws1.Range("C3").Text) + ".txt"
$ws1.Range("E4").Value = ($ws1.Range("C4").Text) + ".txt"
# E5 stays "None" (fail row) - unchanged
$ws1.Range("E6").Value = "testing2.txt"
$ws1.Range("E7").Value = "testing3.txt"
# E8 stays "None" (fail row) - unchanged

# Sheet F2-2
$ws2 = $wb.Worksheets.Item("F2-2")
$ws2.Range("E2").Value = "testing1.php"
# E3, E4, E5 stay "None" (fail rows) - unchanged

# Sheet F2-3
$ws3 = $wb.Worksheets.Item("F2-3")
$ws3.Range("E2").Value = "testing1.php"
# E3, E4, E5 stay "None" (fail rows) - unchanged
$ws3.Range("E6").Value = "testingspecial1.php"
# E7, E8, E9 stay "None" (fail rows) - unchanged

# --- Selection / active-tab bookkeeping ---
# Before: F2-3 tab was active/selected (F9). After: F2-1 is active/selected (E9);
# F2-2's selection moves to E2; F2-3's selection moves to C25 (and is no
# longer the active tab).
$ws2.Range("E2").Select() | Out-Null
$ws3.Range("C25").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("E9").Select() | Out-Null
